$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure target cells keep their original text (inline-string-like) representation
# rather than being auto-converted to numbers/percent values by Excel's type inference.
foreach ($addr in @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "E24", "E25", "D26", "E26", "E27", "E28", "E40", "D41", "E41", "D42", "E42", "E43", "D44", "E44", "D45", "E45", "E46", "E47", "D49", "E49", "D50", "E50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "245.61"
$ws.Range("E2").Value = "-0.52%"
$ws.Range("D3").Value = "30.30"
$ws.Range("E3").Value = "0.67%"
$ws.Range("D4").Value = "5.153"
$ws.Range("E4").Value = "-0.50%"
$ws.Range("D5").Value = "0.05766"
$ws.Range("E5").Value = "0.57%"
$ws.Range("D6").Value = "6.660"
$ws.Range("E6").Value = "0.94%"
$ws.Range("D7").Value = "3.249"
$ws.Range("E7").Value = "5.95%"
$ws.Range("D8").Value = "0.8509"
$ws.Range("E8").Value = "-1.29%"
$ws.Range("D9").Value = "0.8567"
$ws.Range("E9").Value = "-3.06%"
$ws.Range("D10").Value = "0.1384"
$ws.Range("E10").Value = "1.52%"
$ws.Range("D11").Value = "0.07086"
$ws.Range("E11").Value = "-0.03%"
$ws.Range("E12").Value = "13.36%"
$ws.Range("D13").Value = "0.09368"
$ws.Range("E13").Value = "-0.36%"
$ws.Range("D14").Value = "0.001534"
$ws.Range("E14").Value = "1.22%"
$ws.Range("D15").Value = "0.0005938"
$ws.Range("E15").Value = "-0.95%"
$ws.Range("D16").Value = "0.006039"
$ws.Range("E16").Value = "0.98%"
$ws.Range("D17").Value = "3.526"
$ws.Range("E17").Value = "0.60%"
$ws.Range("D18").Value = "2.201"
$ws.Range("E18").Value = "-3.03%"
$ws.Range("D19").Value = "0.3164"
$ws.Range("E19").Value = "-0.64%"
$ws.Range("D20").Value = "0.03386"
$ws.Range("E20").Value = "3.20%"
$ws.Range("E21").Value = "1.21%"
$ws.Range("D22").Value = "3.499"
$ws.Range("E22").Value = "-3.21%"
$ws.Range("D23").Value = "0.04132"
$ws.Range("E23").Value = "-0.39%"
$ws.Range("E24").Value = "2.22%"
$ws.Range("E25").Value = "0.98%"
$ws.Range("D26").Value = "0.004148"
$ws.Range("E26").Value = "-7.95%"
$ws.Range("E27").Value = "-0.80%"
$ws.Range("E28").Value = "4.57%"
$ws.Range("E40").Value = "-0.60%"
$ws.Range("D41").Value = "0.1072"
$ws.Range("E41").Value = "0.04%"
$ws.Range("D42").Value = "0.002469"
$ws.Range("E42").Value = "37.27%"
$ws.Range("E43").Value = "-48.77%"
$ws.Range("D44").Value = "0.008946"
$ws.Range("E44").Value = "-10.92%"
$ws.Range("D45").Value = "0.00005478"
$ws.Range("E45").Value = "7.11%"
$ws.Range("E46").Value = "0.03%"
$ws.Range("E47").Value = "-20.19%"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").Value = "0.03%"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").Value = "0.03%"
